$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "oYbfl297"
$ws.Range("B2").Value = 23081734
$ws.Range("C2").Value = "ldaaikk20"
$ws.Range("D2").Value = "Wh#5j6U&"
$ws.Range("F2").Value = "OrhYOmqy"
$ws.Range("G2").Value = "Rxmn"
